$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A91: was stored as text "11/069" (shared string), should be the
# actual date serial 2020-06-11 (43993), matching the existing date style.
$ws.Range("A91").Value = 43993

# Copy the formatting of the last existing data row down to the new row
# so the new row's styles (date format on A, centered numbers on B:F)
# match the rest of the table instead of defaulting to "General".
$ws.Range("A91:F91").Copy()
$ws.Range("A92:F92").PasteSpecial(-4122) # xlPasteFormats

# Add new row 92 with data for 2020-06-12
$ws.Range("A92").Value = 43994
$ws.Range("B92").Value = 1363
$ws.Range("C92").Value = 463
$ws.Range("D92").Value = 468
$ws.Range("E92").Value = 245
$ws.Range("F92").Value = 48

# Grow the worksheet table so the new row is included.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F92"))

# Mirror the author's final selection/scroll state.
$ws.Range("F92").Select()
